$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 390 (shifts existing rows 390-407 down to 391-408)
$ws.Rows.Item(390).Insert()

# Populate the newly inserted row 390 with the new data record
$ws.Range("A390").Value = 9
$ws.Range("B390").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C390").Value = "Metropolitana"
$ws.Range("D390").Value = 45008
$ws.Range("E390").Value = 13
$ws.Range("F390").Value = 100112043
$ws.Range("G390").Value = "Pepino ensalada"
$ws.Range("H390").Value = "Sin especificar"
$ws.Range("I390").Value = "Primera"
$ws.Range("J390").Value = 70
$ws.Range("K390").Value = 10000
$ws.Range("L390").Value = 11000
$ws.Range("M390").Value = 10500
$ws.Range("N390").Value = "`$/caja 60 unidades"
$ws.Range("O390").Value = "Región de Arica y Parinacota"
$ws.Range("P390").Value = 175
$ws.Range("Q390").Value = 60
$ws.Range("R390").Value = "Hortaliza"
